$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Highlight two more "done-ish" tasks the same way rows 5, 6 and 13 already
# --- are (apply the built-in "Good" cell style). Column C keeps its centered
# --- alignment (column C's base style centers horizontally), so re-apply that
# --- after the style swap since "Good" alone resets alignment to default.
$ws.Range("A7:B7").Style = "Good"
$ws.Range("C7").Style = "Good"
$ws.Range("C7").HorizontalAlignment = -4108

$ws.Range("A10:B10").Style = "Good"
$ws.Range("C10").Style = "Good"
$ws.Range("C10").HorizontalAlignment = -4108

# --- New task row at the bottom of the list.
$ws.Range("A16").Value = "Rūšiavimas skelbimų pagal radimo datą"
$ws.Range("B16").Value = "Ignas"
$ws.Range("C16").Value = "vidutinis"

# --- Move the viewport / selection down to the new area of interest.
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C8").Select() | Out-Null
